$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: fill in the reactivated account data (Xigua) ---
$ws.Range("A15").Value2 = 43749.4583333333
$ws.Range("B15").Value2 = "xxigua@example.com"
$ws.Range("C15").Value2 = "Xavier"
$ws.Range("D15").Value2 = "Xigua"
$ws.Range("E15").Value2 = "Mango"
$ws.Range("F15").NumberFormat = "mm/dd/yy"
$ws.Range("F15").Value2 = 43889
$ws.Range("G15").Value2 = "No"

# New mailto hyperlink on the email cell for row 15
$ws.Hyperlinks.Add($ws.Range("B15"), "mailto:xxigua@example.com", "", "", "xxigua@example.com")

# --- Selection moves from E1 to D16 ---
$ws.Range("D16").Select()

# --- Page setup / margins cosmetic refresh (header/footer margin to 1.3cm, clear explicit first-page-number override) ---
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.HeaderMargin = 36.850393700787386
$ws.PageSetup.FooterMargin = 36.850393700787386
$ws.PageSetup.FirstPageNumber = -4105
